{"js": "// Replace the 25 \"two-digit \u00d7 two-digit = product\" answer strings in the\n// worksheet table with the new values from the commit, one-for-one and in\n// document order. Each old value is unique in the document, so an exact,\n// case-sensitive search reliably targets the single run that holds it.\n// Using search().insertText(..., \"Replace\") edits the text of the existing\n// run in place, leaving paragraph/run formatting (font, size, alignment)\n// untouched, exactly like the diff (which only touches the <w:t> content).\nconst replacements = [\n  [\"85\u00d713=1105\", \"48\u00d774=3552\"],\n  [\"86\u00d764=5504\", \"13\u00d785=1105\"],\n  [\"79\u00d780=6320\", \"47\u00d799=4653\"],\n  [\"33\u00d773=2409\", \"35\u00d766=2310\"],\n  [\"54\u00d791=4914\", \"23\u00d765=1495\"],\n  [\"65\u00d790=5850\", \"14\u00d732=448\"],\n  [\"28\u00d755=1540\", \"78\u00d738=2964\"],\n  [\"63\u00d771=4473\", \"21\u00d750=1050\"],\n  [\"23\u00d790=2070\", \"47\u00d784=3948\"],\n  [\"79\u00d733=2607\", \"28\u00d717=476\"],\n  [\"71\u00d743=3053\", \"45\u00d763=2835\"],\n  [\"80\u00d786=6880\", \"69\u00d749=3381\"],\n  [\"82\u00d793=7626\", \"62\u00d798=6076\"],\n  [\"23\u00d725=575\", \"55\u00d777=4235\"],\n  [\"49\u00d758=2842\", \"28\u00d759=1652\"],\n  [\"76\u00d727=2052\", \"50\u00d736=1800\"],\n  [\"67\u00d764=4288\", \"27\u00d752=1404\"],\n  [\"16\u00d735=560\", \"16\u00d767=1072\"],\n  [\"87\u00d795=8265\", \"43\u00d775=3225\"],\n  [\"84\u00d714=1176\", \"69\u00d735=2415\"],\n  [\"61\u00d718=1098\", \"23\u00d719=437\"],\n  [\"48\u00d761=2928\", \"12\u00d783=996\"],\n  [\"37\u00d758=2146\", \"43\u00d775=3225\"],\n  [\"56\u00d725=1400\", \"84\u00d767=5628\"],\n  [\"80\u00d789=7120\", \"28\u00d729=812\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"two-digit \u00d7 two-digit = product\" answer strings in the\n# worksheet table with the new values from the commit, one-for-one and in\n# document order. Each old value is unique in the document, so a plain\n# Find/Replace (wdReplaceOne) targets exactly the single run that holds it,\n# leaving every other run/paragraph (and all its formatting) untouched -\n# matching the diff, which only rewrites each <w:t> text content.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"85\u00d713=1105\", \"48\u00d774=3552\"),\n    @(\"86\u00d764=5504\", \"13\u00d785=1105\"),\n    @(\"79\u00d780=6320\", \"47\u00d799=4653\"),\n    @(\"33\u00d773=2409\", \"35\u00d766=2310\"),\n    @(\"54\u00d791=4914\", \"23\u00d765=1495\"),\n    @(\"65\u00d790=5850\", \"14\u00d732=448\"),\n    @(\"28\u00d755=1540\", \"78\u00d738=2964\"),\n    @(\"63\u00d771=4473\", \"21\u00d750=1050\"),\n    @(\"23\u00d790=2070\", \"47\u00d784=3948\"),\n    @(\"79\u00d733=2607\", \"28\u00d717=476\"),\n    @(\"71\u00d743=3053\", \"45\u00d763=2835\"),\n    @(\"80\u00d786=6880\", \"69\u00d749=3381\"),\n    @(\"82\u00d793=7626\", \"62\u00d798=6076\"),\n    @(\"23\u00d725=575\",  \"55\u00d777=4235\"),\n    @(\"49\u00d758=2842\", \"28\u00d759=1652\"),\n    @(\"76\u00d727=2052\", \"50\u00d736=1800\"),\n    @(\"67\u00d764=4288\", \"27\u00d752=1404\"),\n    @(\"16\u00d735=560\",  \"16\u00d767=1072\"),\n    @(\"87\u00d795=8265\", \"43\u00d775=3225\"),\n    @(\"84\u00d714=1176\", \"69\u00d735=2415\"),\n    @(\"61\u00d718=1098\", \"23\u00d719=437\"),\n    @(\"48\u00d761=2928\", \"12\u00d783=996\"),\n    @(\"37\u00d758=2146\", \"43\u00d775=3225\"),\n    @(\"56\u00d725=1400\", \"84\u00d767=5628\"),\n    @(\"80\u00d789=7120\", \"28\u00d729=812\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceOne = 1 (replace exactly the first/only match)\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)\n}\n"}
